# The "Förändrad" (Changed) column C date stamp is refreshed for every
# data row in the sheet (row 1 is the header), moving from 2023-10-30
# (serial 45229) to 2023-11-01 (serial 45231).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45231
}
